$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.018.25'
$ws.Range('E2').Value = '  -0.23%  '
$ws.Range('D3').Value = '2.422.55'
$ws.Range('E3').Value = '  +0.05%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '562.99'
$ws.Range('E5').Value = '  -0.18%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.27'
$ws.Range('E6').Value = '  -0.79%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  -0.45%  '
$ws.Range('E9').Value = '  -7.56%  '
$ws.Range('E10').Value = '  -0.21%  '
$ws.Range('E11').Value = '  -0.62%  '
$ws.Range('E12').Value = '  -3.99%  '
$ws.Range('E13').Value = '  -1.34%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.18'
$ws.Range('E14').Value = '  +0.50%  '
$ws.Range('E15').Value = '  -1.83%  '
$ws.Range('E16').Value = '  -0.55%  '
$ws.Range('D17').Value = '61.892.45'
$ws.Range('E17').Value = '  -0.26%  '
$ws.Range('D18').Value = '2.418.37'
$ws.Range('E18').Value = '  +0.02%  '
$ws.Range('E19').Value = '  +0.21%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '323.72'
$ws.Range('E20').Value = '  -0.40%  '
$ws.Range('E21').Value = '  +0.82%  '
$ws.Range('E22').Value = '  -1.28%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  -0.04%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '67.10'
$ws.Range('E24').Value = '  +2.49%  '
$ws.Range('E25').Value = '  +0.78%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.78'
$ws.Range('E26').Value = '  -2.93%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '555.67'
$ws.Range('E27').Value = '  -5.33%  '
$ws.Range('E29').Value = '  -0.20%  '
$ws.Range('D30').Value = '0.0₃0929'
$ws.Range('E30').Value = '  -1.79%  '
$ws.Range('E32').Value = '  -4.41%  '
$ws.Range('E33').Value = '  -2.14%  '
$ws.Range('E34').Value = '  -1.16%  '
$ws.Range('E35').Value = '  -3.22%  '
$ws.Range('E36').Value = '  -0.06%  '
$ws.Range('E37').Value = '  -0.89%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.379'
$ws.Range('E38').Value = '  -1.23%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.46'
$ws.Range('E39').Value = '  -4.70%  '
$ws.Range('B40').Value = 'Monero'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '152.40'
$ws.Range('E40').Value = '  -1.20%  '
$ws.Range('E41').Value = '  -0.34%  '
$ws.Range('E42').Value = '  -1.27%  '
$ws.Range('E43').Value = '  -0.40%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '147.39'
$ws.Range('E44').Value = '  -1.76%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.24'
$ws.Range('E45').Value = '  -4.71%  '
$ws.Range('E46').Value = '  -0.42%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0528'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.596'
$ws.Range('E48').Value = '  +0.52%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '19.83'
$ws.Range('E49').Value = '  -2.83%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0919'
$ws.Range('E50').Value = '  -0.70%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0227'
$ws.Range('E51').Value = '  -0.61%  '
